$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 205.20833
$ws.Range("I6").Value = 205.20833
$ws.Range("K6").Value = 615.62499
$ws.Range("M6").Value = -503.62499

$ws.Range("H17").Value = 1478.4667
$ws.Range("J17").Value = 1478.4667
$ws.Range("L17").Value = 4435.4001
$ws.Range("N17").Value = -4771.4001

$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752

$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142

$ws.Range("H74").Value = 6382.364
$ws.Range("I74").Value = 6240.6
$ws.Range("K74").Value = 6240.6
$ws.Range("M74").Value = -5304.6

$ws.Range("H76").Value = 7342.857
$ws.Range("J76").Value = 7583.3335
$ws.Range("L76").Value = 7583.3335
$ws.Range("N76").Value = -8213.333500000001

$ws.Range("H77").Value = 6382.364
$ws.Range("I77").Value = 6240.6
$ws.Range("K77").Value = 31203
$ws.Range("M77").Value = -26523

$ws.Range("H79").Value = 7342.857
$ws.Range("J79").Value = 7583.3335
$ws.Range("L79").Value = 7583.3335
$ws.Range("N79").Value = -9767.333500000001

$ws.Range("H132").Value = 5160.2104
$ws.Range("I132").Value = 5533.242
$ws.Range("J132").Value = 2698.2
$ws.Range("K132").Value = 16599.726
$ws.Range("L132").Value = 8094.599999999999
$ws.Range("M132").Value = -14069.726
$ws.Range("N132").Value = -13154.6

$ws.Range("H135").Value = 958.1429000000001
$ws.Range("I135").Value = 922.5
$ws.Range("K135").Value = 8302.5
$ws.Range("M135").Value = -5767.5

$ws.Range("H137").Value = 42296.812
$ws.Range("I137").Value = 50253.418
$ws.Range("K137").Value = 150760.254
$ws.Range("M137").Value = -148210.254

$ws.Range("H141").Value = 12397.214
$ws.Range("I141").Value = 18034.334
$ws.Range("K141").Value = 54103.00199999999
$ws.Range("M141").Value = -48923.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4343.317
$ws.Range("I32").Value = 2479.1155
$ws.Range("K32").Value = 2479.1155
$ws.Range("M32").Value = -2192.1155

$ws.Range("H132").Value = 8686.157999999999
$ws.Range("I132").Value = 9216.286
$ws.Range("J132").Value = 7201.8
$ws.Range("K132").Value = 27648.858
$ws.Range("L132").Value = 21605.4
$ws.Range("M132").Value = -25118.858
$ws.Range("N132").Value = -26665.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 32029.428
$ws.Range("J50").Value = 32201.77
$ws.Range("L50").Value = 32201.77
$ws.Range("N50").Value = -33349.77

$ws.Range("H56").Value = 18000
$ws.Range("J56").Value = 18000
$ws.Range("L56").Value = 18000
$ws.Range("N56").Value = -19478

$ws.Range("H82").Value = 37182.59
$ws.Range("I82").Value = 28026
$ws.Range("K82").Value = 28026
$ws.Range("M82").Value = -27643

$ws.Range("H85").Value = 37182.59
$ws.Range("I85").Value = 28026
$ws.Range("K85").Value = 28026
$ws.Range("M85").Value = -26700

$ws.Range("H86").Value = 9994.227999999999
$ws.Range("I86").Value = 16691.727
$ws.Range("J86").Value = 3296.7273
$ws.Range("K86").Value = 16691.727
$ws.Range("L86").Value = 3296.7273
$ws.Range("M86").Value = -15568.727
$ws.Range("N86").Value = -5542.7273

$ws.Range("H89").Value = 9994.227999999999
$ws.Range("I89").Value = 16691.727
$ws.Range("J89").Value = 3296.7273
$ws.Range("K89").Value = 83458.63499999999
$ws.Range("L89").Value = 16483.6365
$ws.Range("M89").Value = -77842.63499999999
$ws.Range("N89").Value = -27715.6365

$ws.Range("H134").Value = 6683.303
$ws.Range("I134").Value = 5658.04
$ws.Range("J134").Value = 9887.25
$ws.Range("K134").Value = 16974.12
$ws.Range("L134").Value = 29661.75
$ws.Range("M134").Value = -14439.12
$ws.Range("N134").Value = -34731.75

$ws.Range("H138").Value = 73556
$ws.Range("J138").Value = 79445
$ws.Range("L138").Value = 79445
$ws.Range("N138").Value = -89725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3552.625
$ws.Range("I58").Value = 3191.9167
$ws.Range("K58").Value = 3191.9167
$ws.Range("M58").Value = -2988.9167

$ws.Range("H132").Value = 40632.19
$ws.Range("I132").Value = 40632.19
$ws.Range("K132").Value = 121896.57
$ws.Range("M132").Value = -119366.57

$ws.Range("H134").Value = 12269.19
$ws.Range("I134").Value = 8693.866
$ws.Range("K134").Value = 26081.598
$ws.Range("M134").Value = -23546.598

$ws.Range("H136").Value = 3552.625
$ws.Range("I136").Value = 3191.9167
$ws.Range("K136").Value = 9575.750100000001
$ws.Range("M136").Value = -7025.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 996.25
$ws.Range("J23").Value = 996.25
$ws.Range("L23").Value = 2988.75
$ws.Range("N23").Value = -3458.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10260.692
$ws.Range("I80").Value = 1898.8
$ws.Range("J80").Value = 15486.875
$ws.Range("K80").Value = 1898.8
$ws.Range("L80").Value = 15486.875
$ws.Range("M80").Value = -900.8
$ws.Range("N80").Value = -17482.875

$ws.Range("H83").Value = 10260.692
$ws.Range("I83").Value = 1898.8
$ws.Range("J83").Value = 15486.875
$ws.Range("K83").Value = 9494
$ws.Range("L83").Value = 77434.375
$ws.Range("M83").Value = -4502
$ws.Range("N83").Value = -87418.375

$ws.Range("H113").Value = 4285.091
$ws.Range("I113").Value = 3427.4
$ws.Range("K113").Value = 3427.4
$ws.Range("M113").Value = -1257.4

$ws.Range("H132").Value = 4304.361
$ws.Range("I132").Value = 2742.1853
$ws.Range("J132").Value = 16355.429
$ws.Range("K132").Value = 8226.555899999999
$ws.Range("L132").Value = 49066.287
$ws.Range("M132").Value = -5696.555899999999
$ws.Range("N132").Value = -54126.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7433.0347
$ws.Range("I40").Value = 5476.6523
$ws.Range("J40").Value = 14932.5
$ws.Range("K40").Value = 5476.6523
$ws.Range("L40").Value = 14932.5
$ws.Range("M40").Value = -5340.6523
$ws.Range("N40").Value = -15204.5

$ws.Range("H68").Value = 5553.875
$ws.Range("J68").Value = 8999
$ws.Range("L68").Value = 8999
$ws.Range("N68").Value = -10497

$ws.Range("H71").Value = 5553.875
$ws.Range("J71").Value = 8999
$ws.Range("L71").Value = 44995
$ws.Range("N71").Value = -52483

$ws.Range("H82").Value = 1925.4445
$ws.Range("I82").Value = 2175.1428
$ws.Range("J82").Value = 1051.5
$ws.Range("K82").Value = 2175.1428
$ws.Range("L82").Value = 1051.5
$ws.Range("M82").Value = -1814.1428
$ws.Range("N82").Value = -1773.5

$ws.Range("H85").Value = 1925.4445
$ws.Range("I85").Value = 2175.1428
$ws.Range("J85").Value = 1051.5
$ws.Range("K85").Value = 2175.1428
$ws.Range("L85").Value = 1051.5
$ws.Range("M85").Value = -927.1428000000001
$ws.Range("N85").Value = -3547.5

$ws.Range("H132").Value = 4000.3
$ws.Range("I132").Value = 3333.8333
$ws.Range("K132").Value = 10001.4999
$ws.Range("M132").Value = -7471.499899999999

$ws.Range("H133").Value = 309999.84
$ws.Range("J133").Value = 309999.84
$ws.Range("L133").Value = 309999.84
$ws.Range("N133").Value = -315059.84

$ws.Range("H136").Value = 26447.582
$ws.Range("I136").Value = 57683.445
$ws.Range("K136").Value = 173050.335
$ws.Range("M136").Value = -170500.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 420.2
$ws.Range("I23").Value = 325
$ws.Range("K23").Value = 325
$ws.Range("M23").Value = -96

$ws.Range("H33").Value = 19666
$ws.Range("I33").Value = 19999.5
$ws.Range("J33").Value = 18999
$ws.Range("K33").Value = 19999.5
$ws.Range("L33").Value = 18999
$ws.Range("M33").Value = -19749.5
$ws.Range("N33").Value = -19499

$ws.Range("H36").Value = 19666
$ws.Range("I36").Value = 19999.5
$ws.Range("J36").Value = 18999
$ws.Range("K36").Value = 19999.5
$ws.Range("L36").Value = 18999
$ws.Range("M36").Value = -19749.5
$ws.Range("N36").Value = -19499

$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2627
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 506471.97
$ws.Range("I132").Value = 9676.538
$ws.Range("K132").Value = 29029.614
$ws.Range("M132").Value = -26499.614

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 7511.778
$ws.Range("I136").Value = 9475.263000000001
$ws.Range("K136").Value = 28425.789
$ws.Range("M136").Value = -25875.789
